$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    3  = -13.119
    21 = -12.63
    23 = -12.68
    25 = -11.988
    53 = -12.575
    57 = -13.742
    59 = -12.788
    69 = -10.613
    79 = -11.999
    83 = -13.012
    93 = -10.281
}

foreach ($row in $updates.Keys) {
    $ws.Range("C$row").Value = $updates[$row]
}
